$wb = $excel.ActiveWorkbook

# Both the "展览" and "全部类型" sheets contain identical data tables and
# both need the same cell updates (matching the diff, which touches the
# same F8/F9/F15/F20 cells in two worksheets).
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F8").Value = 11287
    $ws.Range("F9").Value = 4295
    $ws.Range("F15").Value = 112
    $ws.Range("F20").Value = 11110
}
